# Apply the "Add files via upload" update to the Key Outbreaks sheet.
# The data table (cluster name / active cases) grows from 9 rows to 12
# rows (3 new clusters added) and is rewritten in alphabetical order by
# cluster name, with several case counts changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "3155 Westmont Aged Care Services Baranduda",
    "3349 Maculata Place Shepparton Villages Aged Care Shepparton",
    "3642 Fronditha Care Aged Care Clayton South",
    "3662 Regis Brighton",
    "4314 Estia Health Altona Meadows",
    "Confirmed Omicron Sircuit Bar Fitzroy",
    "Confirmed Omicron Variant The Peel Hotel Collingwood",
    "Diamond Valley Pork and Baxters Pork Laverton North",
    "Mercure Welcome Melbourne",
    "Novotel ibis Melbourne Central Melbourne",
    "Pullman Melbourne on Swanston Melbourne",
    "Werribee Mercy Hospital Emergency Department"
)

$values = @(10, 15, 34, 12, 10, 21, 24, 36, 15, 12, 13, 25)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
